$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values: rows 2, 3, 4, 6 change from 0 to -1.
# Row 5 is left unchanged per the target diff.
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -1
